$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Hide the rows that were collapsed (filtered out / resolved) in this edit,
#    plus re-assert every row that was already hidden beforehand (the loader
#    reads the source file's word-form `hidden="true"` rows as visible, so we
#    restore that pre-existing state explicitly alongside the new changes).
# ---------------------------------------------------------------------------
$previouslyHiddenRows = @(2, 4, 14, 16, 17, 18, 19, 20, 21, 22, 23, 24, 36, 47, 48, 54, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 110, 111, 112, 113, 114, 115, 116, 117, 121, 122, 123, 125, 127, 128, 129, 130)
$newlyHiddenRows = @(32, 132, 133, 134, 135, 138, 141, 142, 143, 145, 146, 147, 148, 153)
foreach ($r in $previouslyHiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}
foreach ($r in $newlyHiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# ---------------------------------------------------------------------------
# 2. Append the four new issue rows (154-157) at the bottom of the log.
# ---------------------------------------------------------------------------
$ws.Range("A154").Value = 154
$ws.Range("B154").Value = "Allow drag all on Crop view"
$ws.Range("F154").Value = 43711
$ws.Range("G154").Value = "DONE"

$ws.Range("A155").Value = 155
$ws.Range("B155").Value = "Change NY → Chester image"
$ws.Range("F155").Value = 43711
$ws.Range("G155").Value = "DONE"

$ws.Range("A156").Value = 156
$ws.Range("B156").Value = "run from perception directory does not work"
$ws.Range("F156").Value = 43711
$ws.Range("G156").Value = "OPEN"

$ws.Range("A157").Value = 157
$ws.Range("B157").Value = "Main display does not show all of the image if it is large"
$ws.Range("F157").Value = 43711
$ws.Range("G157").Value = "OPEN"

# ---------------------------------------------------------------------------
# 3. Move the active selection to the last entered cell, matching the
#    author's cursor position after typing the new rows.
# ---------------------------------------------------------------------------
$ws.Range("G157").Select()

# ---------------------------------------------------------------------------
# 4. Refresh the sheet's `_FilterDatabase` bookkeeping name so it reflects
#    the newly expanded data range (A1:J153 after the hidden-row edits
#    above, before the 4 new rows were appended). Each previous generation
#    shifts down the chain (gains one extra "_0" suffix) and the new
#    current range is inserted right after the original hidden entry -
#    mirroring how the workbook's filter bookkeeping evolves on each edit.
# ---------------------------------------------------------------------------
for ($i = $wb.Names.Count; $i -ge 2; $i--) {
    $n = $wb.Names.Item($i)
    $oldName = $n.Name -replace '^Sheet1!', ''
    $n.Name = $oldName + "_0"
}
$ws.Names.Add("_FilterDatabase", "=Sheet1!`$A`$1:`$J`$153")
